$wb = $excel.ActiveWorkbook

$ws_CUMPLIMIENTOMENSUAL = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws_CUMPLIMIENTOMENSUAL.Range("D7").Value = 444.6
$ws_CUMPLIMIENTOMENSUAL.Range("E7").Value = 1155.4
$ws_CUMPLIMIENTOMENSUAL.Range("F7").Value = 0.277875
$ws_CUMPLIMIENTOMENSUAL.Range("D16").Value = 7261.32
$ws_CUMPLIMIENTOMENSUAL.Range("E16").Value = 14611.78
$ws_CUMPLIMIENTOMENSUAL.Range("F16").Value = 0.3319748915334361
$ws_CUMPLIMIENTOMENSUAL.Range("D18").Value = 2001.17
$ws_CUMPLIMIENTOMENSUAL.Range("E18").Value = -401.1700000000001
$ws_CUMPLIMIENTOMENSUAL.Range("F18").Value = 1.25073125
$ws_CUMPLIMIENTOMENSUAL.Range("D20").Value = 15535.01
$ws_CUMPLIMIENTOMENSUAL.Range("E20").Value = -651.130000000001
$ws_CUMPLIMIENTOMENSUAL.Range("F20").Value = 1.043747329325418
$ws_CUMPLIMIENTOMENSUAL.Range("D21").Value = 193.91
$ws_CUMPLIMIENTOMENSUAL.Range("E21").Value = 452.09
$ws_CUMPLIMIENTOMENSUAL.Range("F21").Value = 0.3001702786377709
$ws_CUMPLIMIENTOMENSUAL.Range("D32").Value = 4993.04
$ws_CUMPLIMIENTOMENSUAL.Range("E32").Value = 15696.96
$ws_CUMPLIMIENTOMENSUAL.Range("F32").Value = 0.2413262445625906
$ws_CUMPLIMIENTOMENSUAL.Range("D33").Value = 14747.8
$ws_CUMPLIMIENTOMENSUAL.Range("E33").Value = 43973.43000000001
$ws_CUMPLIMIENTOMENSUAL.Range("F33").Value = 0.2511493713602388
$ws_CUMPLIMIENTOMENSUAL.Range("D36").Value = 3529.27
$ws_CUMPLIMIENTOMENSUAL.Range("E36").Value = 6441.07304517915
$ws_CUMPLIMIENTOMENSUAL.Range("F36").Value = 0.3539767873590337
$ws_CUMPLIMIENTOMENSUAL.Range("D37").Value = 16832.44
$ws_CUMPLIMIENTOMENSUAL.Range("E37").Value = 10624.5676
$ws_CUMPLIMIENTOMENSUAL.Range("F37").Value = 0.6130471406505347
$ws_CUMPLIMIENTOMENSUAL.Range("D38").Value = 802.67
$ws_CUMPLIMIENTOMENSUAL.Range("E38").Value = 200.33
$ws_CUMPLIMIENTOMENSUAL.Range("F38").Value = 0.8002691924227318
$ws_CUMPLIMIENTOMENSUAL.Range("D41").Value = 1521
$ws_CUMPLIMIENTOMENSUAL.Range("E41").Value = 879
$ws_CUMPLIMIENTOMENSUAL.Range("F41").Value = 0.63375
$ws_CUMPLIMIENTOMENSUAL.Range("D42").Value = 756.66
$ws_CUMPLIMIENTOMENSUAL.Range("E42").Value = 243.34
$ws_CUMPLIMIENTOMENSUAL.Range("F42").Value = 0.75666
$ws_CUMPLIMIENTOMENSUAL.Range("D49").Value = 6163.29
$ws_CUMPLIMIENTOMENSUAL.Range("E49").Value = 7336.71
$ws_CUMPLIMIENTOMENSUAL.Range("F49").Value = 0.45654
$ws_CUMPLIMIENTOMENSUAL.Range("D50").Value = 24512.01
$ws_CUMPLIMIENTOMENSUAL.Range("E50").Value = 31547.69
$ws_CUMPLIMIENTOMENSUAL.Range("F50").Value = 0.4372483263378149
$ws_CUMPLIMIENTOMENSUAL.Range("D52").Value = 6070.05
$ws_CUMPLIMIENTOMENSUAL.Range("E52").Value = -2870.05
$ws_CUMPLIMIENTOMENSUAL.Range("F52").Value = 1.896890625
$ws_CUMPLIMIENTOMENSUAL.Range("D71").Value = 359.21
$ws_CUMPLIMIENTOMENSUAL.Range("E71").Value = 5078.3732
$ws_CUMPLIMIENTOMENSUAL.Range("F71").Value = 0.06606059839231516
$ws_CUMPLIMIENTOMENSUAL.Range("D84").Value = 5747.39
$ws_CUMPLIMIENTOMENSUAL.Range("E84").Value = 44642.78
$ws_CUMPLIMIENTOMENSUAL.Range("F84").Value = 0.1140577616626418
$ws_CUMPLIMIENTOMENSUAL.Range("D104").Value = 71726.64999999999
$ws_CUMPLIMIENTOMENSUAL.Range("E104").Value = -71726.64999999999
$ws_CUMPLIMIENTOMENSUAL.Range("D107").Value = 284.93
$ws_CUMPLIMIENTOMENSUAL.Range("E107").Value = 8383.98
$ws_CUMPLIMIENTOMENSUAL.Range("F107").Value = 0.03286803069820773
$ws_CUMPLIMIENTOMENSUAL.Range("D125").Value = 221033.19
$ws_CUMPLIMIENTOMENSUAL.Range("E125").Value = 281122.5890371132
$ws_CUMPLIMIENTOMENSUAL.Range("F125").Value = 0.4401685676580929
$ws_CUMPLIMIENTOMENSUAL.Columns.Item(6).ColumnWidth = 25.15

$ws_VENTAMENSUAL = $wb.Worksheets.Item("VENTA MENSUAL")
$ws_VENTAMENSUAL.Range("F15").Value = 3799.99
$ws_VENTAMENSUAL.Range("F19").Value = 3172.07
$ws_VENTAMENSUAL.Range("F24").Value = 1599.58
$ws_VENTAMENSUAL.Range("F70").Value = 3757.42
$ws_VENTAMENSUAL.Range("F78").Value = 1606.1
$ws_VENTAMENSUAL.Range("F91").Value = 373.31
$ws_VENTAMENSUAL.Range("F94").Value = 1541.99
$ws_VENTAMENSUAL.Range("F103").Value = 8065.41
$ws_VENTAMENSUAL.Range("F107").Value = 307.93
$ws_VENTAMENSUAL.Range("F110").Value = 409.57
$ws_VENTAMENSUAL.Range("F112").Value = 6346.89
$ws_VENTAMENSUAL.Range("F120").Value = 5209.92
$ws_VENTAMENSUAL.Range("F124").Value = 15847.07
$ws_VENTAMENSUAL.Range("F138").Value = 322.11
$ws_VENTAMENSUAL.Range("F142").Value = 2063.26
$ws_VENTAMENSUAL.Range("F194").Value = 290.97
$ws_VENTAMENSUAL.Range("F211").Value = 91.58
$ws_VENTAMENSUAL.Range("F212").Value = 183.09
$ws_VENTAMENSUAL.Range("F224").Value = 1190.95
$ws_VENTAMENSUAL.Range("F280").Value = 366.34
$ws_VENTAMENSUAL.Range("F304").Value = 201329.79

$ws_VENTASPORGRUPO = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws_VENTASPORGRUPO.Range("M15").Value = 1895.34
$ws_VENTASPORGRUPO.Range("H19").Value = 444.6
$ws_VENTASPORGRUPO.Range("M19").Value = 2021.43
$ws_VENTASPORGRUPO.Range("O24").Value = 1599.58
$ws_VENTASPORGRUPO.Range("D70").Value = 933.12
$ws_VENTASPORGRUPO.Range("E70").Value = 193.91
$ws_VENTASPORGRUPO.Range("M70").Value = 1229.62
$ws_VENTASPORGRUPO.Range("M78").Value = 1606.1
$ws_VENTASPORGRUPO.Range("E91").Value = 69.13
$ws_VENTASPORGRUPO.Range("D94").Value = 468.29
$ws_VENTASPORGRUPO.Range("H94").Value = 1073.7
$ws_VENTASPORGRUPO.Range("M103").Value = 6512.36
$ws_VENTASPORGRUPO.Range("M107").Value = 307.93
$ws_VENTASPORGRUPO.Range("I110").Value = 26.1
$ws_VENTASPORGRUPO.Range("M110").Value = 383.47
$ws_VENTASPORGRUPO.Range("D112").Value = 2289.6
$ws_VENTASPORGRUPO.Range("C120").Value = 1555.2
$ws_VENTASPORGRUPO.Range("D120").Value = 2799.36
$ws_VENTASPORGRUPO.Range("L120").Value = 855.36
$ws_VENTASPORGRUPO.Range("O124").Value = 4136.14
$ws_VENTASPORGRUPO.Range("E138").Value = 142.11
$ws_VENTASPORGRUPO.Range("H138").Value = 180
$ws_VENTASPORGRUPO.Range("D142").Value = 570.24
$ws_VENTASPORGRUPO.Range("E142").Value = 55.65
$ws_VENTASPORGRUPO.Range("M142").Value = 777.39
$ws_VENTASPORGRUPO.Range("M194").Value = 143.97
$ws_VENTASPORGRUPO.Range("D211").Value = 91.58
$ws_VENTASPORGRUPO.Range("D212").Value = 88.53
$ws_VENTASPORGRUPO.Range("M224").Value = 1190.95
$ws_VENTASPORGRUPO.Range("D280").Value = 366.34
$ws_VENTASPORGRUPO.Range("C304").Value = "7 de 302"
$ws_VENTASPORGRUPO.Range("D304").Value = "26 de 302"
$ws_VENTASPORGRUPO.Range("E304").Value = "7 de 302"
$ws_VENTASPORGRUPO.Range("H304").Value = "8 de 302"
$ws_VENTASPORGRUPO.Range("I304").Value = "8 de 302"
$ws_VENTASPORGRUPO.Range("L304").Value = "19 de 302"
$ws_VENTASPORGRUPO.Range("M304").Value = "39 de 302"
$ws_VENTASPORGRUPO.Range("O304").Value = "9 de 302"
